$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1) with login info
$ws.Range("C1").Value = "user1"
$ws.Range("E1").Value = "password1"

# Row 4: C4 course changes from CPSC 3415 to CYBR 4125, D4 credits 1 -> 3
$ws.Range("C4").Value = "CYBR 4125"
$ws.Range("D4").Value = 3

# Row 5: A5 course changes from GEOL 1110 to STAT 1401, C5 changes from CYBR 4125 to CPSC 4135
$ws.Range("A5").Value = "STAT 1401"
$ws.Range("C5").Value = "CPSC 4135"

# Row 6: C6 changes from CPSC 4135 to CYBR 4145
$ws.Range("C6").Value = "CYBR 4145"

# New row 15: CPSC 4205, 3 credits
$ws.Range("A15").Value = "CPSC 4205"
$ws.Range("B15").Value = 3
